$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Draw": remove the two oldest fixtures (25-12-2024) and add the new
# Iran fixture at the bottom.
# ---------------------------------------------------------------------------
$wsDraw = $wb.Worksheets.Item("Draw")
$wsDraw.Range("A2:A3").EntireRow.Delete()

$wsDraw.Range("A9").Value = "27-12-2024 13:00"
$wsDraw.Range("B9").Value = "IRAN"
$wsDraw.Range("C9").Value = "PERSIAN GULF PRO LEAGUE"
$wsDraw.Range("D9").Value = "Shams Azar Qazvin - Gol Gohar"
$wsDraw.Range("E9").Value = 66.7
$wsDraw.Range("F9").Value = 2.8

# ---------------------------------------------------------------------------
# Sheet "Btts": append the Philippines - Thailand fixture.
# ---------------------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")
$wsBtts.Range("A12").Value = "27-12-2024 13:00"
$wsBtts.Range("B12").Value = "WORLD"
$wsBtts.Range("C12").Value = "AFF CHAMPIONSHIP"
$wsBtts.Range("D12").Value = "Philippines - Thailand"
$wsBtts.Range("E12").Value = 76.7
$wsBtts.Range("F12").Value = 1.95

# ---------------------------------------------------------------------------
# Sheet "Over_Under": append the three new Israeli fixtures.
# ---------------------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")

$wsOU.Range("A7").Value = "27-12-2024 12:00"
$wsOU.Range("B7").Value = "ISRAEL"
$wsOU.Range("C7").Value = "LIGA ALEF"
$wsOU.Range("D7").Value = "Maccabi K. Ata Bialik - Maccabi Nujeidat Ahmed"
$wsOU.Range("E7").Value = 80
$wsOU.Range("F7").Value = 1.6
$wsOU.Range("G7").Value = 65
$wsOU.Range("H7").Value = 2.55

$wsOU.Range("A8").Value = "27-12-2024 12:00"
$wsOU.Range("B8").Value = "ISRAEL"
$wsOU.Range("C8").Value = "STATE CUP"
$wsOU.Range("D8").Value = "Hapoel Afula - Ashdod"
$wsOU.Range("E8").Value = 75
$wsOU.Range("F8").Value = 1.9
$wsOU.Range("G8").Value = 70
$wsOU.Range("H8").Value = 3.3

$wsOU.Range("A9").Value = "27-12-2024 12:00"
$wsOU.Range("B9").Value = "ISRAEL"
$wsOU.Range("C9").Value = "STATE CUP"
$wsOU.Range("D9").Value = "Hapoel Ramat HaSharon - Maccabi Herzliya"
$wsOU.Range("E9").Value = 76.3
$wsOU.Range("F9").Value = 1.8
$wsOU.Range("G9").Value = 61.3
$wsOU.Range("H9").Value = 3

# ---------------------------------------------------------------------------
# Sheet "Away Win": append three new fixtures.
# ---------------------------------------------------------------------------
$wsAway = $wb.Worksheets.Item("Away Win")

$wsAway.Range("A4").Value = "26-12-2024 17:30"
$wsAway.Range("B4").Value = "ISRAEL"
$wsAway.Range("C4").Value = "STATE CUP"
$wsAway.Range("D4").Value = "Bnei Sakhnin - Ironi Modi'in"
$wsAway.Range("E4").Value = 70
$wsAway.Range("F4").Value = 9

$wsAway.Range("A5").Value = "26-12-2024 15:00"
$wsAway.Range("B5").Value = "WORLD"
$wsAway.Range("C5").Value = "AFRICAN NATIONS CHAMPIONSHIP - QUALIFICATION"
$wsAway.Range("D5").Value = "Burundi - Uganda"
$wsAway.Range("E5").Value = 73.3
$wsAway.Range("F5").Value = 2.15

$wsAway.Range("A6").Value = "27-12-2024 14:00"
$wsAway.Range("B6").Value = "INDIA"
$wsAway.Range("C6").Value = "INDIAN SUPER LEAGUE"
$wsAway.Range("D6").Value = "Mohammedan - Odisha"
$wsAway.Range("E6").Value = 80
$wsAway.Range("F6").Value = 1.85
